$d = $word.ActiveDocument

# The 4th table in the document (the "don vi thi cong" / contractor block,
# directly below the "Cong ty TNHH Nghia Van" line) switches from an
# auto-fit width to a fixed dxa width, and its 3rd column grows from
# 3645 dxa (182.25 pt) to 3933 dxa (196.65 pt). Word COM widths are
# expressed in points, so dxa values are divided by 20.

$tbl = $d.Tables.Item(4)

$tbl.PreferredWidthType = 3          # dxa (fixed width)
$tbl.PreferredWidth = 432.55         # 8651 dxa = 425 + 4293 + 3933

$col = $tbl.Columns.Item(3)
$col.PreferredWidthType = 3          # dxa
$col.Width = 196.65                  # 3933 dxa
